{"js": "// Fill in the \"Burning Glass\" self-analysis paragraph (for the author,\n// \"Connor\") and drop the stray _GoBack bookmark left over from the last\n// edit session.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate \"Natalie\" (the last of the named headings before the two blank\n// paragraphs that need to become the new analysis paragraph).\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet natalieIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Natalie\") {\n    natalieIdx = i;\n    break;\n  }\n}\nif (natalieIdx === -1) {\n  throw new Error(\"Could not find the 'Natalie' heading paragraph\");\n}\n\nconst firstBlank = paragraphs.items[natalieIdx + 1];\nconst secondBlank = paragraphs.items[natalieIdx + 2];\n\nconst analysisText =\n  \"After analysing the data, it\u2019s clear to me that the job of a UI Game Programmer\" +\n  \" is quite a niche title in the IT world. UX designer (user experience) is the only\" +\n  \" similar alternative, and is ranked in the bottom five job titles. This has made me\" +\n  \" evaluate what draws me to the job in order to branch out into other titles. One of\" +\n  \" the key things that attracts me to it is design, and seeing that graphic design and\" +\n  \" creativity are still highly sought-after skills, I think I would there would be\" +\n  \" opportunities in other fields that combine design and code skills.\";\n\n// Put the new text into the first of the two empty paragraphs, then remove\n// the now-redundant second empty paragraph entirely.\nfirstBlank.insertText(analysisText, \"Start\");\nsecondBlank.delete();\nawait context.sync();\n\n// Remove the leftover \"_GoBack\" bookmark (Word drops this automatically on\n// a normal save; the diff shows it gone from the final paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Fill in the \"Burning Glass\" self-analysis paragraph (for the author,\n# \"Connor\") and drop the stray _GoBack bookmark left over from the last\n# edit session.\n\n$d = $word.ActiveDocument\n\n# Find the \"Natalie\" heading paragraph; the new analysis text replaces the\n# first of the two blank paragraphs that follow it, and the second blank\n# paragraph is removed outright.\n$natalieIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq \"Natalie\") {\n        $natalieIdx = $i\n        break\n    }\n}\nif ($natalieIdx -eq -1) {\n    throw \"Could not find the 'Natalie' heading paragraph\"\n}\n\n$analysisText = \"After analysing the data, it\u2019s clear to me that the job of a UI Game Programmer is quite a niche title in the IT world. UX designer (user experience) is the only similar alternative, and is ranked in the bottom five job titles. This has made me evaluate what draws me to the job in order to branch out into other titles. One of the key things that attracts me to it is design, and seeing that graphic design and creativity are still highly sought-after skills, I think I would there would be opportunities in other fields that combine design and code skills.\"\n\n$firstBlank = $d.Paragraphs.Item($natalieIdx + 1)\n$firstBlank.Range.Text = $analysisText\n\n$secondBlank = $d.Paragraphs.Item($natalieIdx + 2)\n$secondBlank.Range.Delete()\n\n# Remove the leftover \"_GoBack\" bookmark (Word drops this automatically on\n# a normal save; the diff shows it gone from the final paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
